$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D11").Value = 16
$ws.Range("D13").Value = 0.5
$ws.Range("D14").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("D19").Value = 0.5
$ws.Range("D21").Value = 0.03

$ws.Range("D11").Select()
